$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (rows 2 and 3), pushing the
# existing data rows (old 2-5) down to become rows 4-7.
$ws.Rows("2:3").Insert()

# The insert copies formatting from the row above (the bold header row),
# so clear that back to the plain/default formatting used by the rest of
# the data rows.
$ws.Range("A2:T3").ClearFormats()

# Restore the date number format on column D for the two new rows, matching
# the other data rows (style index with numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 2: new "Primera" quality entry for the 2021-12-02 date ---
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44532
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100101
$ws.Cells.Item(2, 8).Value = "Berries"
$ws.Cells.Item(2, 9).Value = 100101004
$ws.Cells.Item(2, 10).Value = "Frambuesa"
$ws.Cells.Item(2, 11).Value = "Sin especificar"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 100
$ws.Cells.Item(2, 14).Value = 10000
$ws.Cells.Item(2, 15).Value = 10000
$ws.Cells.Item(2, 16).Value = 10000
$ws.Cells.Item(2, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(2, 18).Value = "Región de Ñuble"
$ws.Cells.Item(2, 19).Value = 5000
$ws.Cells.Item(2, 20).Value = 2

# --- Row 3: new "Segunda" quality entry for the 2021-12-02 date ---
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44532
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101004
$ws.Cells.Item(3, 10).Value = "Frambuesa"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Segunda"
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 8000
$ws.Cells.Item(3, 15).Value = 8000
$ws.Cells.Item(3, 16).Value = 8000
$ws.Cells.Item(3, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(3, 18).Value = "Región de Ñuble"
$ws.Cells.Item(3, 19).Value = 4000
$ws.Cells.Item(3, 20).Value = 2
